$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new "Product Icon" value for row 2 (column G), using a new shared string
$ws.Range("G2").Value = "D:\Work Files\Graphics\Test Storefront\Test product image.png"

# Update "Operation" column (A) for rows 4 and 5 from "Change" to "Skip"
$ws.Range("A4").Value = "Skip"
$ws.Range("A5").Value = "Skip"

# Update the active selection on the bottom-right (frozen) pane
$ws.Activate()
$ws.Range("G3").Select()
